$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new week's label and contribution percentages
$ws.Range("B3").Value = "Jan 27-Feb3rd"
$ws.Range("C3").Value = 0.25
$ws.Range("D3").Value = 0.15
$ws.Range("E3").Value = 0.25
$ws.Range("F3").Value = 0.15
$ws.Range("G3").Value = 0.1
$ws.Range("H3").Value = 0.1

# Update the active selection
$ws.Range("G6").Select()
